$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, shifting existing rows 67-79 down to 68-80
$ws.Rows.Item(67).Insert()

# Populate the new row 67 with values (copy of row, then adjust the new record's data)
$ws.Cells.Item(67, 1).Value = 8
$ws.Cells.Item(67, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(67, 3).Value = "Coquimbo"
$ws.Cells.Item(67, 4).Value = 44641
$ws.Cells.Item(67, 5).Value = 4
$ws.Cells.Item(67, 6).Value = 100112030
$ws.Cells.Item(67, 7).Value = "Poroto granado"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 500
$ws.Cells.Item(67, 11).Value = 28000
$ws.Cells.Item(67, 12).Value = 29000
$ws.Cells.Item(67, 13).Value = 28500
$ws.Cells.Item(67, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(67, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(67, 16).Value = 1140
$ws.Cells.Item(67, 17).Value = 25
$ws.Cells.Item(67, 18).Value = "Hortaliza"
